$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the GitHub-Actions-refreshed crypto price/volume snapshot.
# Values are stored as literal text in this sheet (prices like "28.047.23" and
# percentages like "  +3.60%  " are display strings, not numbers), so for any
# new value that Excel would otherwise auto-parse as a number we force the
# cell to text format first and then assign the literal string.

$ws.Range('D2').Value = '28.074.46'
$ws.Range('E2').Value = '  +3.72%  '
$ws.Range('D3').Value = '1.727.08'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.63'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('E6').Value = '  +1.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.07'
$ws.Range('E8').Value = '  +13.51%  '
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0903'
$ws.Range('E11').Value = '  +2.22%  '
$ws.Range('D12').Value = '1.971.96'
$ws.Range('E12').Value = '  +3.12%  '
$ws.Range('D13').Value = '1.725.99'
$ws.Range('E13').Value = '  +3.09%  '
$ws.Range('E14').Value = '  +3.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.563'
$ws.Range('E15').Value = '  +5.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.72'
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('D17').Value = '28.016.72'
$ws.Range('E17').Value = '  +3.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.47'
$ws.Range('E18').Value = '  +2.50%  '
$ws.Range('D19').Value = '0.0₃0758'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.88'
$ws.Range('E20').Value = '  -3.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.64'
$ws.Range('E22').Value = '  +3.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.77'
$ws.Range('E23').Value = '  +4.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.13'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.15'
$ws.Range('E25').Value = '  +1.81%  '
$ws.Range('E26').Value = '  +4.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.78'
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('E28').Value = '  +2.10%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  +2.63%  '
$ws.Range('E31').Value = '  +2.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.44'
$ws.Range('E32').Value = '  +2.91%  '
$ws.Range('E33').Value = '  +2.82%  '
$ws.Range('D34').Value = '1.491.70'
$ws.Range('E34').Value = '  -3.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.66'
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.959'
$ws.Range('E36').Value = '  +3.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.610'
$ws.Range('E37').Value = '  +1.73%  '
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('E39').Value = '  +1.11%  '
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '71.14'
$ws.Range('E41').Value = '  +5.29%  '
$ws.Range('E42').Value = '  +4.33%  '
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('D45').Value = '1.875.71'
$ws.Range('E45').Value = '  +2.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.796'
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.76'
$ws.Range('E47').Value = '  +12.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '91.62'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').Value = '0.0₆0110'
$ws.Range('E49').Value = '  +2.15%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.105'
$ws.Range('E50').Value = '  +1.15%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.24'
$ws.Range('E51').Value = '  +2.44%  '
